$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "SkillCode"
$ws.Range("B1").Value = "Skill Description"
$ws.Range("C1").Value = "SFIA Level"
$ws.Range("D1").Value = "Keycode"
$ws.Range("E1").Value = "Description"

$ws.Range("A2").Value = "Autonomy"
$ws.Range("B2").Value = "Autonomy"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Works under close direction"

$ws.Range("A3").Value = "Autonomy"
$ws.Range("B3").Value = "Autonomy"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Uses little discretion in attending to enquiries"

$ws.Range("A4").Value = "Autonomy"
$ws.Range("B4").Value = "Autonomy"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "Is expected to seek guidance in unexpected situations"

$ws.Range("A5").Value = "Influence"
$ws.Range("B5").Value = "Influence"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Minimal influence"

$ws.Range("A6").Value = "Influence"
$ws.Range("B6").Value = "Influence"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "May work alone or interact with immediate colleagues"

$ws.Range("A7").Value = "Complexity"
$ws.Range("B7").Value = "Complexity"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Performs routine activities in a structured environment"

$ws.Range("A8").Value = "Complexity"
$ws.Range("B8").Value = "Complexity"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Requires assistance in resolving unexpected problems"

$ws.Range("A9").Value = "Complexity"
$ws.Range("B9").Value = "Complexity"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = "Participates in the generation of new ideas"

$ws.Range("A10").Value = "Knowledge"
$ws.Range("B10").Value = "Knowledge"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "Has a basic generic knowledge appropriate to area of work"

$ws.Range("A11").Value = "Knowledge"
$ws.Range("B11").Value = "Knowledge"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Applies newly acquired knowledge to develop new skills"

$ws.Range("A12").Value = "SINT"
$ws.Range("B12").Value = "Systems integration and build"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = "Produces software builds from software source code"

$ws.Range("A13").Value = "SINT"
$ws.Range("B13").Value = "Systems integration and build"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "Conducts tests as defined in an integration test specification and records the details of any failures"

$ws.Range("A14").Value = "SINT"
$ws.Range("B14").Value = "Systems integration and build"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "Analyses and reports on integration test activities and results"

$ws.Range("A15").Value = "SINT"
$ws.Range("B15").Value = "Systems integration and build"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = "Identifies and reports issues and risks"

$ws.Range("A16").Value = "TEST"
$ws.Range("B16").Value = "Testing"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "Designs test cases, creates test scripts and test data, and automates repeatable tasks working to the requirements or specifications provided"

$ws.Range("A17").Value = "TEST"
$ws.Range("B17").Value = "Testing"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = "Defines test conditions for given requirements"

$ws.Range("A18").Value = "TEST"
$ws.Range("B18").Value = "Testing"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = "Executes and records manual and automated testing  in accordance with test plans"

$ws.Range("A19").Value = "TEST"
$ws.Range("B19").Value = "Testing"
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = "Analyses and reports on test activities, results, issues and risks"

$ws.Range("A20").Value = "ITOP"
$ws.Range("B20").Value = "IT infrastructure"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = "Carries out routine operational procedures, including the execution of specified automation tools/scripts"

$ws.Range("A21").Value = "ITOP"
$ws.Range("B21").Value = "IT infrastructure"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = "Amends existing automation tasks under supervision to gain a basic understanding of the scripting language/automation tools"

$ws.Range("A22").Value = "ITOP"
$ws.Range("B22").Value = "IT infrastructure"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = "Contributes to maintenance and installation"

$ws.Range("A23").Value = "ITOP"
$ws.Range("B23").Value = "IT infrastructure"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = "Monitors and reports on infrastructure performance to enable service delivery"

$ws.Range("A24").Value = "ITOP"
$ws.Range("B24").Value = "IT infrastructure"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = "Resolves issues or refers to others for assistance"

$ws.Range("A25").Value = "NTAS"
$ws.Range("B25").Value = "Network support"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = "Contributes to the operational configuration of network components"

$ws.Range("A26").Value = "NTAS"
$ws.Range("B26").Value = "Network support"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = "Assists in the investigation and resolution of network problems"

$ws.Range("A27").Value = "NTAS"
$ws.Range("B27").Value = "Network support"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = "Assists with specified maintenance procedures"

$ws.Range("A28").Value = "HSIN"
$ws.Range("B28").Value = "Systems installation and removal"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = "Installs or removes system components using supplied installation instructions and tools"

$ws.Range("A29").Value = "HSIN"
$ws.Range("B29").Value = "Systems installation and removal"
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = "Conducts standard tests and contributes to investigations of problems and faults"

$ws.Range("A30").Value = "HSIN"
$ws.Range("B30").Value = "Systems installation and removal"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = "Confirms the correct working of installations"

$ws.Range("A31").Value = "HSIN"
$ws.Range("B31").Value = "Systems installation and removal"
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = "Documents results in accordance with agreed procedures"

$ws.Range("A32").Value = "MADE"
$ws.Range("B32").Value = "MADE"
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = "Example MT"

$ws.Range("A34").Value = "CFMG"
$ws.Range("B34").Value = "Configuration management"
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = "Applies tools, techniques and processes to administer, track, log, report on and correct configuration items, components and changes"

$ws.Range("A35").Value = "CFMG"
$ws.Range("B35").Value = "Configuration management"
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = "Assists with audits to check the accuracy of the information and undertakes any necessary corrective action under direction"
